$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row labels (bold, bordered, centered) - B1:G1
$ws.Range("B1").Value = "IT.CEL.SETS.P2:LIE"
$ws.Range("C1").Value = "IT.MLT.MAIN.P2:LIE"
$ws.Range("D1").Value = "IT.NET.USER.P2:LIE"
$ws.Range("E1").Value = "SG.GEN.PARL.ZS:LIE"
$ws.Range("F1").Value = "SP.RUR.TOTL:LIE"
$ws.Range("G1").Value = "SP.URB.TOTL:LIE"

# Row labels column A (bold, bordered, centered) - A2:A13
$ws.Range("A2").Value = "IT.MLT.MAIN.P2:LIE:cor-value"
$ws.Range("A3").Value = "IT.MLT.MAIN.P2:LIE:p-value"
$ws.Range("A4").Value = "IT.NET.USER.P2:LIE:cor-value"
$ws.Range("A5").Value = "IT.NET.USER.P2:LIE:p-value"
$ws.Range("A6").Value = "SG.GEN.PARL.ZS:LIE:cor-value"
$ws.Range("A7").Value = "SG.GEN.PARL.ZS:LIE:p-value"
$ws.Range("A8").Value = "SP.POP.TOTL:LIE:cor-value"
$ws.Range("A9").Value = "SP.POP.TOTL:LIE:p-value"
$ws.Range("A10").Value = "SP.RUR.TOTL:LIE:cor-value"
$ws.Range("A11").Value = "SP.RUR.TOTL:LIE:p-value"
$ws.Range("A12").Value = "SP.URB.TOTL:LIE:cor-value"
$ws.Range("A13").Value = "SP.URB.TOTL:LIE:p-value"

# Numeric correlation / p-value data
$ws.Range("B2").Value = -0.9124778156941027
$ws.Range("D2").Value = -0.9440610758751974
$ws.Range("B3").Value = 0.000005361934690974171
$ws.Range("D3").Value = 0.000000391867107818568
$ws.Range("B4").Value = 0.9281316918743813
$ws.Range("B5").Value = 0.00000170166428896518
$ws.Range("B6").Value = 0.3673964864245902
$ws.Range("C6").Value = -0.110966381142038
$ws.Range("D6").Value = 0.1629344763150769
$ws.Range("F6").Value = 0.259638478163536
$ws.Range("G6").Value = 0.1791428937464186
$ws.Range("B7").Value = 0.1962619069242975
$ws.Range("C7").Value = 0.7056850054987932
$ws.Range("D7").Value = 0.5778367970409857
$ws.Range("F7").Value = 0.370029736600021
$ws.Range("G7").Value = 0.5400165598975177
$ws.Range("B8").Value = 0.9689759135900431
$ws.Range("C8").Value = -0.9720388250918793
$ws.Range("D8").Value = 0.9808984567071966
$ws.Range("E8").Value = 0.2533683558565097
$ws.Range("F8").Value = 0.9999216904422149
$ws.Range("G8").Value = 0.9897901683085839
$ws.Range("B9").Value = 0.00000001204028696122471
$ws.Range("C9").Value = 0.000000006496196113023744
$ws.Range("D9").Value = 0.0000000006730684946563117
$ws.Range("E9").Value = 0.3821072968163085
$ws.Range("F9").Value = 0.000000000000000000000003328953645971784
$ws.Range("G9").Value = 0.00000000001599863050636106
$ws.Range("B10").Value = 0.9700462451543288
$ws.Range("C10").Value = -0.9694590740747172
$ws.Range("D10").Value = 0.9817306939862932
$ws.Range("G10").Value = 0.9879289377404633
$ws.Range("B11").Value = 0.00000000977594151986932
$ws.Range("C11").Value = 0.0000000109696094728465
$ws.Range("D11").Value = 0.0000000005161323129250377
$ws.Range("G11").Value = 0.0000000000435215668501791
$ws.Range("B12").Value = 0.946028223712011
$ws.Range("C12").Value = -0.9906292103237435
$ws.Range("D12").Value = 0.9605301458546136
$ws.Range("B13").Value = 0.0000003174865896649551
$ws.Range("C13").Value = 0.000000000009581032331929928
$ws.Range("D13").Value = 0.00000005012632987627068

# Apply header/label styling: bold font, thin border, centered horizontal, top vertical.
# Build the style on a single anchor cell first (keeps the style table minimal), then
# copy just the formatting onto the rest of the header row / label column.
$styleAnchor = $ws.Range("B1")
$styleAnchor.Font.Bold = $true
$styleAnchor.HorizontalAlignment = -4108
$styleAnchor.VerticalAlignment = -4160
$styleAnchor.Borders.LineStyle = 1
$styleAnchor.Borders.Weight = 2
$styleAnchor.Copy()
$ws.Range("C1:G1").PasteSpecial(-4122)
$ws.Range("A2:A13").PasteSpecial(-4122)
$excel.CutCopyMode = $false
